$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Strapping" pin-function column (old column L) is being relocated so it
# sits right after "Common Restrictions" (column E), ahead of "ADC" - i.e. it
# becomes the new column F. Cutting column L and inserting the cut cells
# before column F shifts the old F:K block right into G:L while leaving the
# following Sensor/Camera columns (M/N) untouched.
$ws.Columns("L:L").Cut() | Out-Null
$ws.Columns("F:F").Insert() | Out-Null

# Re-home three pin functions so the design stops relying on GPIO14 (an ADC2
# pin, which conflicts with wifi's use of ADC2): the "Audio in from SLIC via
# op-amp" signal moves off pin 14 (row 18) onto pin 32 (row 36, was "RM to
# SLIC"), which in turn moves onto pin 22 (row 26); "FR to SLIC" moves from
# pin 33 (row 37) onto pin 23 (row 27).
$pin14 = $ws.Range("C18").Value2
$pin32 = $ws.Range("C36").Value2
$pin33 = $ws.Range("C37").Value2

$ws.Range("C36").Value = $pin14
$ws.Range("C26").Value = $pin32
$ws.Range("C27").Value = $pin33
$ws.Range("C18").ClearContents() | Out-Null
$ws.Range("C37").ClearContents() | Out-Null

# The legend entry describing pins with boot restrictions now also calls out
# the new wifi/ADC2 restriction.
$ws.Range("C47").Value = "pin has boot and/or wifi restrictions"
